# Updates cryptos list snapshot (price + 1h volume change columns),
# and fixes the ordering of the Aptos / InternetComputer(DFINITY) rows.
# Columns: A=Rank(unused) B=Coin C=Link D=Price E=Volume(1h)
# Numeric-looking Price strings are assigned with a leading apostrophe
# so Excel keeps them as text (preserving formatting such as trailing
# zeros) instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '62.911.88'
$ws.Cells.Item(2, 5).Value = '  -0.17%  '
$ws.Cells.Item(3, 4).Value = '2.581.73'
$ws.Cells.Item(3, 5).Value = '  +1.19%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).Value = '''583.33'
$ws.Cells.Item(5, 5).Value = '  +0.71%  '
$ws.Cells.Item(6, 4).Value = '''146.55'
$ws.Cells.Item(6, 5).Value = '  -0.32%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 4).Value = '''0.594'
$ws.Cells.Item(8, 5).Value = '  +2.01%  '
$ws.Cells.Item(9, 5).Value = '  +2.46%  '
$ws.Cells.Item(10, 5).Value = '  +2.38%  '
$ws.Cells.Item(11, 5).Value = '  -0.18%  '
$ws.Cells.Item(12, 5).Value = '  -0.27%  '
$ws.Cells.Item(13, 4).Value = '''27.24'
$ws.Cells.Item(13, 5).Value = '  +0.20%  '
$ws.Cells.Item(14, 4).Value = '3.043.84'
$ws.Cells.Item(14, 5).Value = '  +1.24%  '
$ws.Cells.Item(15, 4).Value = '62.760.35'
$ws.Cells.Item(15, 5).Value = '  -0.29%  '
$ws.Cells.Item(16, 5).Value = '  +3.18%  '
$ws.Cells.Item(17, 4).Value = '2.579.63'
$ws.Cells.Item(17, 5).Value = '  +1.02%  '
$ws.Cells.Item(18, 4).Value = '''11.30'
$ws.Cells.Item(18, 5).Value = '  -0.28%  '
$ws.Cells.Item(19, 4).Value = '''341.08'
$ws.Cells.Item(19, 5).Value = '  +1.71%  '
$ws.Cells.Item(20, 4).Value = '''4.37'
$ws.Cells.Item(20, 5).Value = '  +0.89%  '
$ws.Cells.Item(21, 4).Value = '''6.68'
$ws.Cells.Item(21, 5).Value = '  -0.84%  '
$ws.Cells.Item(22, 4).Value = '''1.00'
$ws.Cells.Item(22, 5).Value = '  +0.00%  '
$ws.Cells.Item(23, 4).Value = '''5.71'
$ws.Cells.Item(23, 5).Value = '  -0.70%  '
$ws.Cells.Item(24, 4).Value = '''67.01'
$ws.Cells.Item(24, 5).Value = '  +2.43%  '
$ws.Cells.Item(25, 4).Value = '2.709.19'
$ws.Cells.Item(25, 5).Value = '  +1.09%  '
$ws.Cells.Item(26, 5).Value = '  -1.56%  '
$ws.Cells.Item(27, 5).Value = '  -0.80%  '
$ws.Cells.Item(28, 4).Value = '''0.999'
$ws.Cells.Item(28, 5).Value = '  -0.11%  '
$ws.Cells.Item(29, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(29, 4).Value = '''8.33'
$ws.Cells.Item(29, 5).Value = '  +0.03%  '
$ws.Cells.Item(30, 2).Value = 'Aptos'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(30, 4).Value = '''7.83'
$ws.Cells.Item(30, 5).Value = '  +6.85%  '
$ws.Cells.Item(31, 5).Value = '  -2.78%  '
$ws.Cells.Item(32, 4).Value = '''1.93'
$ws.Cells.Item(32, 5).Value = '  +2.24%  '
$ws.Cells.Item(33, 4).Value = '0.0₃0820'
$ws.Cells.Item(33, 5).Value = '  +0.89%  '
$ws.Cells.Item(34, 4).Value = '''464.53'
$ws.Cells.Item(34, 5).Value = '  +13.83%  '
$ws.Cells.Item(35, 4).Value = '''175.01'
$ws.Cells.Item(35, 5).Value = '  -1.58%  '
$ws.Cells.Item(36, 5).Value = '  +3.78%  '
$ws.Cells.Item(37, 5).Value = '  +0.10%  '
$ws.Cells.Item(38, 4).Value = '''0.400'
$ws.Cells.Item(38, 5).Value = '  +0.30%  '
$ws.Cells.Item(39, 4).Value = '''18.96'
$ws.Cells.Item(39, 5).Value = '  -0.76%  '
$ws.Cells.Item(40, 5).Value = '  +3.92%  '
$ws.Cells.Item(41, 5).Value = '  +0.03%  '
$ws.Cells.Item(42, 4).Value = '''1.71'
$ws.Cells.Item(42, 5).Value = '  -1.67%  '
$ws.Cells.Item(43, 4).Value = '''157.98'
$ws.Cells.Item(44, 5).Value = '  -0.02%  '
$ws.Cells.Item(45, 4).Value = '''0.635'
$ws.Cells.Item(45, 5).Value = '  +5.49%  '
$ws.Cells.Item(46, 4).Value = '''21.07'
$ws.Cells.Item(46, 5).Value = '  +1.14%  '
$ws.Cells.Item(47, 5).Value = '  +0.58%  '
$ws.Cells.Item(48, 5).Value = '  +0.03%  '
$ws.Cells.Item(49, 5).Value = '  -0.53%  '
$ws.Cells.Item(50, 4).Value = '''18.40'
$ws.Cells.Item(50, 5).Value = '  +1.16%  '
$ws.Cells.Item(51, 4).Value = '''1.72'
$ws.Cells.Item(51, 5).Value = '  +0.16%  '
